$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text format
# so Excel stores them as strings (matching the source data) instead of
# converting them to numeric values.
$textCells = @('D5', 'D8', 'D10', 'D11', 'D13', 'D14', 'D15', 'D19', 'D20', 'D22', 'D25', 'D27', 'D28', 'D29', 'D30', 'D34', 'D39', 'D43', 'D44', 'D45', 'D48', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.286.08'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.044.93'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '228.61'
$ws.Range('E5').Value = '  -2.27%  '
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '55.98'
$ws.Range('E8').Value = '  -4.43%  '
$ws.Range('E9').Value = '  -2.84%  '
$ws.Range('D10').Value = '0.0815'
$ws.Range('E10').Value = '  +3.85%  '
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('D12').Value = '2.345.44'
$ws.Range('D13').Value = '14.56'
$ws.Range('E13').Value = '  -3.84%  '
$ws.Range('D14').Value = '20.54'
$ws.Range('E14').Value = '  -3.59%  '
$ws.Range('D15').Value = '0.752'
$ws.Range('E15').Value = '  -3.67%  '
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '2.041.35'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '37.195.20'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '6.00'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').Value = '69.68'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').Value = '225.62'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -5.53%  '
$ws.Range('E26').Value = '  -2.89%  '
$ws.Range('D27').Value = '168.16'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -5.40%  '
$ws.Range('D29').Value = '1.39'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = '18.93'
$ws.Range('E30').Value = '  -3.02%  '
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('E33').Value = '  -3.68%  '
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  -3.16%  '
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D39').Value = '5.36'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').Value = '1.499.83'
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('E41').Value = '  -7.17%  '
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '95.71'
$ws.Range('E43').Value = '  -6.39%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '16.70'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '0.0933'
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('E47').Value = '  -4.93%  '
$ws.Range('D48').Value = '7.12'
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').Value = '2.230.63'
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').Value = '3.68'
$ws.Range('E51').Value = '  -11.38%  '
